$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new forecast-vintage column BB (one quarter ahead of the existing
# last column BA) and one new trailing row (83) for the newest quarter.
# ---------------------------------------------------------------------------

# Row 1 header: new vintage date in BB1 (copy BA1's date format, then set value)
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# Rows 2-71: BB is an exact copy (value + formatting) of BA for those rows -
# the forecast for these older quarters hasn't changed with the new vintage.
for ($r = 2; $r -le 71; $r++) {
    $ws.Range("BA$r").Copy($ws.Range("BB$r"))
}

# Rows 72-82: the newest quarters get revised forecast values in column BB
# (column BA is left untouched). Formatting is still copied from BA so any
# (here, default) cell formatting is preserved consistently.
$bbNewValues = @{
    72 = "-0.5"
    73 = "-0.3"
    74 = "0.1510250314585848"
    75 = "-0.002717645765269422"
    76 = "-0.136875628108055"
    77 = "-0.05727680561923214"
    78 = "-0.02351950773205924"
    79 = "-0.05824996978169377"
    80 = "-0.0640752710254897"
    81 = "-0.05053803152041143"
    82 = "-0.05090228640356263"
}
for ($r = 72; $r -le 82; $r++) {
    $ws.Range("BA$r").Copy($ws.Range("BB$r"))
    $ws.Range("BB$r").Value = $bbNewValues[$r]
}

# Row 83: brand-new trailing row for the newest quarter - only the date (A)
# and the new vintage forecast (BB) are populated.
$ws.Range("A82").Copy($ws.Range("A83"))
$ws.Range("A83").Value = 46934

$ws.Range("BA82").Copy($ws.Range("BB83"))
$ws.Range("BB83").Value = -0.05570051844454853
